$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), matching the formatting of the other
# header cells (bold, centered, bordered) by copying G1's format.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" column values for each data row.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
